# Rename the worksheet tab / sheet name.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Lista de asistencia Evento 1"

# Remove the second data row (DIEGO BELTRAN LOPEZ entry) entirely,
# shifting the dimension back down to A1:F1.
$ws.Rows.Item(2).Delete()
